$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "AB2" = 16.5
    "AC2" = 12
    "AD2" = 20
    "AJ2" = 19
    "AK2" = 15
    "AM2" = 50
    "AN2" = 5.6
    "AO2" = 29
    "F2" = 1.67
    "G2" = 1.69
    "H2" = 4.9
    "I2" = 5.1
    "J2" = 4.7
    "K2" = 4.9
    "L2" = 1.24
    "M2" = 1.02
    "N2" = 7.8
    "O2" = 1.13
    "P2" = 3.25
    "Q2" = 1.43
    "R2" = 1.91
    "S2" = 2.06
    "T2" = 1.51
    "U2" = 2.9
    "V2" = 1.24
    "W2" = 2.44
    "X2" = 40
    "Y2" = 34
    "Z2" = 50
    "AA3" = 80
    "AB3" = 13
    "AC3" = 9.199999999999999
    "AD3" = 16
    "AE3" = 42
    "AF3" = 14
    "AG3" = 10.5
    "AI3" = 42
    "AM3" = 65
    "AO3" = 32
    "F3" = 1.94
    "G3" = 1.96
    "H3" = 4.1
    "I3" = 4.2
    "L3" = 1.31
    "N3" = 5.3
    "O3" = 1.21
    "P3" = 2.48
    "Q3" = 1.64
    "R3" = 1.59
    "S3" = 2.58
    "T3" = 1.62
    "U3" = 2.48
    "V3" = 1.31
    "X3" = 22
    "Y3" = 21
    "F4" = 1.81
    "G4" = 1.89
    "H4" = 3.9
    "J4" = 4.2
    "O4" = 1.2
    "V4" = 1.28
    "W4" = 2.12
    "X4" = 26
    "F5" = 2.16
    "G5" = 2.74
    "H5" = 3.25
    "I5" = 4.5
    "N5" = 2.4
    "O5" = 1.01
    "Q5" = 2.16
    "V5" = 1.29
    "W5" = 1.57
    "X5" = 19
    "F6" = 1.51
    "G6" = 1.64
    "I6" = 12
    "K6" = 6.4
    "V6" = 1.09
    "W6" = 2.56
    "X6" = 1000
    "F7" = 12.5
    "I7" = 1.26
    "K7" = 8.199999999999999
    "H8" = 1.35
    "J8" = 3.05
    "N8" = 1.1
    "AF9" = 13.5
    "M9" = 1.08
    "U9" = 2.06
    "X9" = 12.5
    "AN10" = 6.6
    "AO10" = 85
    "F10" = 1.55
    "G10" = 1.56
    "M10" = 1.04
    "Q10" = 1.67
    "V10" = 1.16
    "W10" = 2.78
    "X10" = 24
    "Z10" = 60
    "O11" = 1.17
    "Q11" = 1.55
    "S11" = 2.32
    "T11" = 1.55
    "AA12" = 55
    "AB12" = 11
    "AC12" = 8.199999999999999
    "AD12" = 13.5
    "AE12" = 36
    "AI12" = 44
    "AJ12" = 32
    "AK12" = 24
    "AN12" = 18
    "F12" = 2.34
    "G12" = 2.36
    "H12" = 3.2
    "I12" = 3.25
    "K12" = 3.8
    "V12" = 1.44
    "W12" = 1.74
    "Y12" = 13
    "Z12" = 22
    "AB13" = 30
    "AO13" = 5.8
    "L13" = 1.29
    "O13" = 1.23
    "Q13" = 1.69
    "T13" = 1.97
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}

Write-Output "Done applying $($changes.Count) cell changes"